# Update "想去人数" (interest count) values in column F on sheets "展览" and
# "全部类型" to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Row -> (old, new) value for sheet "展览"
$sheetExhibition = $wb.Worksheets.Item("展览")
$exhibitionChanges = @{
    5  = 704
    6  = 129
    9  = 2639
    11 = 1652
    15 = 846
    16 = 121
    20 = 42
    22 = 5819
    24 = 1078
    25 = 122
    29 = 236
    31 = 1069
    32 = 847
}

foreach ($row in $exhibitionChanges.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionChanges[$row]
}

# Row -> new value for sheet "全部类型"
$sheetAll = $wb.Worksheets.Item("全部类型")
$allChanges = @{
    7  = 704
    8  = 129
    14 = 2639
    16 = 1652
    21 = 846
    22 = 121
    25 = 42
    27 = 5819
    29 = 1078
    30 = 122
    34 = 236
    36 = 1069
    37 = 847
}

foreach ($row in $allChanges.Keys) {
    $sheetAll.Cells.Item($row, 6).Value = $allChanges[$row]
}
